$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "1.003") are not converted to floating point numbers by Excel,
# matching the original inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.697.47'
$ws.Range("D3").Value = '1.754.21'
$ws.Range("D4").Value = '1.003'
$ws.Range("D5").Value = '324.46'
$ws.Range("D6").Value = '0.9990'
$ws.Range("D7").Value = '0.4293'
$ws.Range("D8").Value = '0.3642'
$ws.Range("D10").Value = '0.07489'
$ws.Range("D11").Value = '1.121'
$ws.Range("D12").Value = '1.000'
$ws.Range("D13").Value = '21.65'
$ws.Range("D14").Value = '6.153'
$ws.Range("D15").Value = '7.259'
$ws.Range("D16").Value = '1.747.13'
$ws.Range("D17").Value = '0.00001069'
$ws.Range("D18").Value = '87.96'
$ws.Range("D19").Value = '0.06204'
$ws.Range("D20").Value = '0.9994'
$ws.Range("D23").Value = '0.5264'
$ws.Range("D24").Value = '27.720.73'
$ws.Range("D25").Value = '11.69'
$ws.Range("D27").Value = '20.55'
$ws.Range("D28").Value = '152.77'
$ws.Range("D29").Value = '2.372'
$ws.Range("D30").Value = '1.947.91'
$ws.Range("D31").Value = '1.227'
$ws.Range("D32").Value = '127.38'
$ws.Range("D33").Value = '5.739'
$ws.Range("D34").Value = '0.09157'
$ws.Range("D35").Value = '3.657'
$ws.Range("D36").Value = '12.73'
$ws.Range("D37").Value = '0.02315'
$ws.Range("D38").Value = '0.2157'
$ws.Range("D39").Value = '5.123'
$ws.Range("D40").Value = '0.6493'
$ws.Range("D41").Value = '0.06110'
$ws.Range("D42").Value = '1.197'
$ws.Range("D43").Value = '1.429'
$ws.Range("D44").Value = '7.982'
$ws.Range("D45").Value = '0.9988'
$ws.Range("D46").Value = '13.80'
$ws.Range("D47").Value = '0.5948'
$ws.Range("D48").Value = '3.751'
$ws.Range("D49").Value = '126.24'
$ws.Range("D50").Value = '1.974'
$ws.Range("D51").Value = '0.06911'

# Restore default (unstyled) cell style now that the values are locked in as text
$ws.Range("D2:D51").Style = "Normal"

$ws.Range("E2").Value = '  -2.07%  '
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E5").Value = '  -4.15%  '
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  -7.68%  '
$ws.Range("E8").Value = '  -4.27%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  -3.17%  '
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("E13").Value = '  -3.62%  '
$ws.Range("E14").Value = '  -3.04%  '
$ws.Range("E15").Value = '  -3.39%  '
$ws.Range("E16").Value = '  -3.43%  '
$ws.Range("E17").Value = '  -2.38%  '
$ws.Range("E18").Value = '  +7.99%  '
$ws.Range("E19").Value = '  -7.86%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("E22").Value = '  -4.24%  '
$ws.Range("E23").Value = '  -6.01%  '
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("E26").Value = '  -3.89%  '
$ws.Range("E27").Value = '  -0.57%  '
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -3.36%  '
$ws.Range("E31").Value = '  -1.76%  '
$ws.Range("E32").Value = '  -4.00%  '
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("E34").Value = '  -4.83%  '
$ws.Range("E35").Value = '  -9.47%  '
$ws.Range("E36").Value = '  +5.43%  '
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("E38").Value = '  -7.51%  '
$ws.Range("E39").Value = '  -3.02%  '
$ws.Range("E40").Value = '  -1.96%  '
$ws.Range("E41").Value = '  -3.90%  '
$ws.Range("E42").Value = '  -4.02%  '
$ws.Range("E43").Value = '  -4.17%  '
$ws.Range("E44").Value = '  -4.59%  '
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E46").Value = '  -2.71%  '
$ws.Range("E47").Value = '  -3.12%  '
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("E49").Value = '  -3.71%  '
$ws.Range("E50").Value = '  -3.46%  '
$ws.Range("E51").Value = '  -3.42%  '

